$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure target cells keep their original text (inline-string) representation
# rather than being auto-converted to numeric/percentage values by Excel.
$targetCells = @("D2","E2","D3","E3","D4","E4","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","E20","E21","D22","E22","D23","E23","E24","D25","E25","E26","E27","D40","E40","D41","E41","D42","E42","E43","D44","E44","D45","E45","D46","E46","D48","E48","D49","E49","D50","E50")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "261.23"
$ws.Range("E2").Value = "0.09%"
$ws.Range("D3").Value = "26.91"
$ws.Range("E3").Value = "-1.71%"
$ws.Range("D4").Value = "4.716"
$ws.Range("E4").Value = "0.17%"
$ws.Range("D6").Value = "6.729"
$ws.Range("E6").Value = "0.76%"
$ws.Range("D7").Value = "0.8492"
$ws.Range("E7").Value = "0.40%"
$ws.Range("D8").Value = "0.9117"
$ws.Range("E8").Value = "-1.26%"
$ws.Range("D9").Value = "0.1402"
$ws.Range("E9").Value = "-0.12%"
$ws.Range("D10").Value = "0.04961"
$ws.Range("E10").Value = "0.34%"
$ws.Range("D11").Value = "0.07086"
$ws.Range("E11").Value = "-0.37%"
$ws.Range("D12").Value = "0.03114"
$ws.Range("E12").Value = "0.55%"
$ws.Range("D13").Value = "0.09056"
$ws.Range("E13").Value = "-0.16%"
$ws.Range("D14").Value = "0.001532"
$ws.Range("E14").Value = "0.02%"
$ws.Range("D15").Value = "0.0006174"
$ws.Range("E15").Value = "1.72%"
$ws.Range("D16").Value = "0.005977"
$ws.Range("E16").Value = "-2.43%"
$ws.Range("D17").Value = "3.448"
$ws.Range("E17").Value = "-0.08%"
$ws.Range("D18").Value = "3.174"
$ws.Range("E18").Value = "0.99%"
$ws.Range("E19").Value = "-0.85%"
$ws.Range("E20").Value = "-0.38%"
$ws.Range("E21").Value = "1.76%"
$ws.Range("D22").Value = "4.109"
$ws.Range("E22").Value = "0.49%"
$ws.Range("D23").Value = "0.04246"
$ws.Range("E23").Value = "0.33%"
$ws.Range("E24").Value = "-3.22%"
$ws.Range("D25").Value = "0.004071"
$ws.Range("E25").Value = "4.09%"
$ws.Range("E26").Value = "0.00%"
$ws.Range("E27").Value = "4.07%"
$ws.Range("D40").Value = "0.03938"
$ws.Range("E40").Value = "1.83%"
$ws.Range("D41").Value = "0.1113"
$ws.Range("E41").Value = "0.02%"
$ws.Range("D42").Value = "0.004136"
$ws.Range("E42").Value = "0.98%"
$ws.Range("E43").Value = "-2.69%"
$ws.Range("D44").Value = "0.01319"
$ws.Range("E44").Value = "-19.33%"
$ws.Range("D45").Value = "0.00005163"
$ws.Range("E45").Value = "0.32%"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "0.00%"
$ws.Range("D48").Value = "0.2491"
$ws.Range("E48").Value = "84.13%"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.00%"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.00%"
